$d = $word.ActiveDocument

# The journal table (first/only table in the document) lists entries from
# most-recent to oldest. A new entry for "15 FÉVRIER" is added at the top;
# every other entry keeps its own original content and simply shifts down
# by one row.
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add($t.Rows.Item(1))

$dateCell = $newRow.Cells.Item(1)
$dateCell.Range.Text = "15 FÉVRIER"

$contentCell = $newRow.Cells.Item(2)
$bullet1 = "Mise à jour du contrôleur Cie (changement des noms de fonctions reliées au model docs et les alertes)"
$bullet2 = "Ajout de la vue Info de la compagnie"
$bullet3 = "Modifications dans les vues de l’intern soient log, list et index"
$contentCell.Range.Text = $bullet1 + "`r" + $bullet2 + "`r" + $bullet3
